$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.022103
$ws.Range("H2").Value = 0.06630900000000001
$ws.Range("I2").Value = 0.0007043476645371027
$ws.Range("J2").Value = 0.0007043476645371028
$ws.Range("M2").Value = 0.6692693333333334
$ws.Range("N2").Value = 2.007808
$ws.Range("O2").Value = 0.004126561180566838
$ws.Range("P2").Value = 0.004126561180566839
$ws.Range("Q2").Value = 0.01479286007466667
$ws.Range("R2").Value = 0.133135740672
$ws.Range("S2").Value = "0.000002906533730101722091032548"
$ws.Range("T2").Value = "0.000002906533730101722938065495"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.022103
$ws.Range("H3").Value = 0.06630900000000001
$ws.Range("I3").Value = 0.0007043476645371027
$ws.Range("J3").Value = 0.0007043476645371028
$ws.Range("O3").Value = 0.9916964991825307
$ws.Range("P3").Value = 0.9916964991825309
$ws.Range("Q3").Value = 3.555024851692333
$ws.Range("R3").Value = 31.995223665231
$ws.Range("S3").Value = 0.0006984991131288363
$ws.Range("T3").Value = 0.0006984991131288365
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.022103
$ws.Range("H4").Value = 0.06630900000000001
$ws.Range("I4").Value = 0.0007043476645371027
$ws.Range("J4").Value = 0.0007043476645371028
$ws.Range("M4").Value = 0.5637343333333333
$ws.Range("N4").Value = 1.691203
$ws.Range("O4").Value = 0.003475856580040611
$ws.Range("P4").Value = 0.003475856580040611
$ws.Range("Q4").Value = 0.01246021996966667
$ws.Range("R4").Value = 0.112141979727
$ws.Range("S4").Value = "0.000002448211464417525191526129"
$ws.Range("T4").Value = "0.000002448211464417526038559077"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.022103
$ws.Range("H5").Value = 0.06630900000000001
$ws.Range("I5").Value = 0.0007043476645371027
$ws.Range("J5").Value = 0.0007043476645371028
$ws.Range("M5").Value = 0.1137056666666667
$ws.Range("N5").Value = 0.341117
$ws.Range("O5").Value = 0.0007010830568617209
$ws.Range("P5").Value = 0.0007010830568617211
$ws.Range("Q5").Value = 0.002513236350333333
$ws.Range("R5").Value = 0.022619127153
$ws.Range("S5").Value = "0.000000493806213747085838671998"
$ws.Range("T5").Value = "0.000000493806213747086050430235"
$ws.Range("G6").Value = 5.827140333333332
$ws.Range("I6").Value = 0.1856912041222136
$ws.Range("J6").Value = 0.1856912041222136
$ws.Range("M6").Value = 0.6692693333333334
$ws.Range("N6").Value = 2.007808
$ws.Range("O6").Value = 0.004126561180566838
$ws.Range("P6").Value = 0.004126561180566839
$ws.Range("Q6").Value = 3.899926326129778
$ws.Range("R6").Value = 35.099336935168
$ws.Range("S6").Value = 0.0007662661145034394
$ws.Range("T6").Value = 0.0007662661145034395
$ws.Range("G7").Value = 5.827140333333332
$ws.Range("I7").Value = 0.1856912041222136
$ws.Range("J7").Value = 0.1856912041222136
$ws.Range("O7").Value = 0.9916964991825307
$ws.Range("P7").Value = 0.9916964991825309
$ws.Range("Q7").Value = 937.2315386734263
$ws.Range("R7").Value = 8435.083848060836
$ws.Range("S7").Value = 0.1841493170569879
$ws.Range("T7").Value = 0.1841493170569879
$ws.Range("G8").Value = 5.827140333333332
$ws.Range("I8").Value = 0.1856912041222136
$ws.Range("J8").Value = 0.1856912041222136
$ws.Range("M8").Value = 0.5637343333333333
$ws.Range("N8").Value = 1.691203
$ws.Range("O8").Value = 0.003475856580040611
$ws.Range("P8").Value = 0.003475856580040611
$ws.Range("Q8").Value = 3.284959071051444
$ws.Range("R8").Value = 29.56463163946299
$ws.Range("S8").Value = 0.0006454359937038603
$ws.Range("T8").Value = 0.0006454359937038603
$ws.Range("G9").Value = 5.827140333333332
$ws.Range("I9").Value = 0.1856912041222136
$ws.Range("J9").Value = 0.1856912041222136
$ws.Range("M9").Value = 0.1137056666666667
$ws.Range("N9").Value = 0.341117
$ws.Range("O9").Value = 0.0007010830568617209
$ws.Range("P9").Value = 0.0007010830568617211
$ws.Range("Q9").Value = 0.6625788763618887
$ws.Range("R9").Value = 5.963209887256999
$ws.Range("S9").Value = 0.0001301849570183353
$ws.Range("T9").Value = 0.0001301849570183353
$ws.Range("G10").Value = 15.496839
$ws.Range("H10").Value = 46.490517
$ws.Range("I10").Value = 0.4938317132225258
$ws.Range("J10").Value = 0.4938317132225258
$ws.Range("M10").Value = 0.6692693333333334
$ws.Range("N10").Value = 2.007808
$ws.Range("O10").Value = 0.004126561180566838
$ws.Range("P10").Value = 0.004126561180566839
$ws.Range("Q10").Value = 10.371559106304
$ws.Range("R10").Value = 93.344031956736
$ws.Range("S10").Value = 0.002037826777516891
$ws.Range("T10").Value = 0.002037826777516891
$ws.Range("G11").Value = 15.496839
$ws.Range("H11").Value = 46.490517
$ws.Range("I11").Value = 0.4938317132225258
$ws.Range("J11").Value = 0.4938317132225258
$ws.Range("O11").Value = 0.9916964991825307
$ws.Range("P11").Value = 0.9916964991825309
$ws.Range("Q11").Value = 2492.496392692166
$ws.Range("R11").Value = 22432.4675342295
$ws.Range("S11").Value = 0.4897311811880903
$ws.Range("T11").Value = 0.4897311811880904
$ws.Range("G12").Value = 15.496839
$ws.Range("H12").Value = 46.490517
$ws.Range("I12").Value = 0.4938317132225258
$ws.Range("J12").Value = 0.4938317132225258
$ws.Range("M12").Value = 0.5637343333333333
$ws.Range("N12").Value = 1.691203
$ws.Range("O12").Value = 0.003475856580040611
$ws.Range("P12").Value = 0.003475856580040611
$ws.Range("Q12").Value = 8.736100202438999
$ws.Range("R12").Value = 78.624901821951
$ws.Range("S12").Value = 0.001716488209837244
$ws.Range("T12").Value = 0.001716488209837245
$ws.Range("G13").Value = 15.496839
$ws.Range("H13").Value = 46.490517
$ws.Range("I13").Value = 0.4938317132225258
$ws.Range("J13").Value = 0.4938317132225258
$ws.Range("M13").Value = 0.1137056666666667
$ws.Range("N13").Value = 0.341117
$ws.Range("O13").Value = 0.0007010830568617209
$ws.Range("P13").Value = 0.0007010830568617211
$ws.Range("Q13").Value = 1.762078409721
$ws.Range("R13").Value = 15.858705687489
$ws.Range("S13").Value = 0.0003462170470813091
$ws.Range("T13").Value = 0.0003462170470813092
$ws.Range("G14").Value = 10.03472733333333
$ws.Range("H14").Value = 30.104182
$ws.Range("I14").Value = 0.3197727349907235
$ws.Range("J14").Value = 0.3197727349907235
$ws.Range("M14").Value = 0.6692693333333334
$ws.Range("N14").Value = 2.007808
$ws.Range("O14").Value = 0.004126561180566838
$ws.Range("P14").Value = 0.004126561180566839
$ws.Range("Q14").Value = 6.715935272561779
$ws.Range("R14").Value = 60.44341745305601
$ws.Range("S14").Value = 0.001319561754816407
$ws.Range("T14").Value = 0.001319561754816407
$ws.Range("G15").Value = 10.03472733333333
$ws.Range("H15").Value = 30.104182
$ws.Range("I15").Value = 0.3197727349907235
$ws.Range("J15").Value = 0.3197727349907235
$ws.Range("O15").Value = 0.9916964991825307
$ws.Range("P15").Value = 0.9916964991825309
$ws.Range("Q15").Value = 1613.975706915638
$ws.Range("R15").Value = 14525.78136224074
$ws.Range("S15").Value = 0.3171175018243236
$ws.Range("T15").Value = 0.3171175018243237
$ws.Range("G16").Value = 10.03472733333333
$ws.Range("H16").Value = 30.104182
$ws.Range("I16").Value = 0.3197727349907235
$ws.Range("J16").Value = 0.3197727349907235
$ws.Range("M16").Value = 0.5637343333333333
$ws.Range("N16").Value = 1.691203
$ws.Range("O16").Value = 0.003475856580040611
$ws.Range("P16").Value = 0.003475856580040611
$ws.Range("Q16").Value = 5.656920323438445
$ws.Range("R16").Value = 50.912282910946
$ws.Range("S16").Value = 0.001111484165035089
$ws.Range("T16").Value = 0.001111484165035089
$ws.Range("G17").Value = 10.03472733333333
$ws.Range("H17").Value = 30.104182
$ws.Range("I17").Value = 0.3197727349907235
$ws.Range("J17").Value = 0.3197727349907235
$ws.Range("M17").Value = 0.1137056666666667
$ws.Range("N17").Value = 0.341117
$ws.Range("O17").Value = 0.0007010830568617209
$ws.Range("P17").Value = 0.0007010830568617211
$ws.Range("Q17").Value = 1.141005361254889
$ws.Range("R17").Value = 10.269048251294
$ws.Range("S17").Value = 0.0002241872465483294
$ws.Range("T17").Value = 0.0002241872465483295
